# Update YEAR_START and YEAR_END values on the STANDARD_DEFINITION sheet
# YEAR_START (column C): 1900 -> 1700
# YEAR_END   (column D): 2020 -> 2040

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STANDARD_DEFINITION")

$ws.Range("C2:C4").Value = 1700
$ws.Range("D2:D4").Value = 2040
